$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin prices in column D are stored as plain text in this sheet (so values
# using "." as a thousands separator, e.g. "68.718.02", and values with
# trailing zeros, e.g. "1.00", render exactly as scraped instead of being
# normalised as numbers). For any new price whose text would otherwise be
# auto-parsed by Excel as a number, the cell is pre-formatted as Text so the
# literal string is preserved; values that already fail numeric parsing (e.g.
# "68.718.02" has two dots) do not need this and are left with their original
# formatting.

$ws.Range("D2").Value = "68.718.02"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "2.710.94"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.06"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.19"
$ws.Range("E6").Value = "  +3.50%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("D9").Value = "2.708.89"
$ws.Range("E9").Value = "  +2.28%  "
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.33"
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.362"
$ws.Range("E13").Value = "  +2.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.43"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("D15").Value = "3.208.42"
$ws.Range("E15").Value = "  +2.45%  "
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "68.716.04"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "2.754.95"
$ws.Range("E18").Value = "  +4.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.88"
$ws.Range("E19").Value = "  +4.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.70"
$ws.Range("E20").Value = "  +4.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "365.62"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.55"
$ws.Range("E22").Value = "  +2.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.92"
$ws.Range("E23").Value = "  +2.23%  "
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("E25").Value = "  -1.76%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +1.86%  "
$ws.Range("D28").Value = "2.842.54"
$ws.Range("E28").Value = "  +2.35%  "
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "595.64"
$ws.Range("E30").Value = "  +6.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.30"
$ws.Range("E32").Value = "  +2.83%  "
$ws.Range("E33").Value = "  +3.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.96"
$ws.Range("E34").Value = "  +5.01%  "
$ws.Range("E35").Value = "  +2.89%  "
$ws.Range("E36").Value = "  +5.36%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.91"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "160.64"
$ws.Range("E39").Value = "  +0.57%  "
$ws.Range("E40").Value = "  +2.30%  "
$ws.Range("E41").Value = "  +2.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.43"
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("E43").Value = "  +2.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.00"
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "158.10"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("E48").Value = "  +5.21%  "
$ws.Range("E49").Value = "  +5.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.608"
$ws.Range("E50").Value = "  +7.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.09"
$ws.Range("E51").Value = "  -0.49%  "

# Row 45/46 swap: BabyDogeCoin (45) <-> USDe (46), with updated D/E values
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0316"
$ws.Range("E46").Value = "  -5.65%  "
